# Generate Report for Handback
# - Overview sheet: status message updated for both locale columns (E2, F2)
# - zh-cn / de-de sheets: record the generated handback target + handback
#   file + handback datetime, with a hyperlink on the "Latest Target File"
#   cell pointing back at the source markdown file (same target/display as
#   the existing A2 hyperlink).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: handback status -------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.md",
    $null,
    $null,
    "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.md"
) | Out-Null
$zhcn.Range("J2").Value = "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.4ab74e330ad8179519427b2cff08588ae293166e.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 11:07:38"

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Hyperlinks.Add(
    $dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f33054fa241d7766091054952f872c849148cc50/e2e/4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.md",
    $null,
    $null,
    "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.md"
) | Out-Null
$dede.Range("J2").Value = "4c08c1d8-8ae1-4579-8b63-a5087eb7f6dc.4ab74e330ad8179519427b2cff08588ae293166e.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 11:07:44"
